# Applies the PoolingContract_results.xlsx update described by the commit
# "test all the examples": refreshed solver benchmark numbers on the "DE"
# sheet, plus two new result blocks (Solver=Antigone, Solver=SCIP), and the
# relocation of the "Number of nodes" summary table further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DE")

# ---------------------------------------------------------------------
# Update the existing "Solver=BARON" block (rows 2-4) with refreshed
# benchmark numbers.
# ---------------------------------------------------------------------
$ws.Range("C2").Value = -1338.24
$ws.Range("D2").Value = -1338.24

$ws.Range("C3").Value = -1337.2
$ws.Range("D3").Value = -1338.25
$ws.Range("F3").Value = 3005

$ws.Range("C4").Value = -1221.32
$ws.Range("D4").Value = -1338.25
$ws.Range("F4").Value = 10000

# ---------------------------------------------------------------------
# New block: "Solver=Antigone" (header row 6, data rows 7-9)
# ---------------------------------------------------------------------
$ws.Range("A6").Value = "Solver=Antigone"
$ws.Range("B6").Value = "number of scenario"
$ws.Range("C6").Value = "UB"
$ws.Range("D6").Value = "LB"
$ws.Range("E6").Value = "gap"
$ws.Range("F6").Value = "Wall Time"

$ws.Range("B7").Value = 3
$ws.Range("C7").Value = -1338.2470000000001
$ws.Range("D7").Value = -1339.585
$ws.Range("C7").NumberFormat = "0.00"
$ws.Range("D7").NumberFormat = "0.00"
$ws.Range("E7").Formula = "=-(C7-D7)/D7"
$ws.Range("F7").Value = 16

$ws.Range("B8").Value = 9
$ws.Range("C8").Value = -1338.2470000000001
$ws.Range("D8").Value = -1339.585
$ws.Range("C8").NumberFormat = "0.00"
$ws.Range("D8").NumberFormat = "0.00"
$ws.Range("E8").Formula = "=-(C8-D8)/D8"
$ws.Range("F8").Value = 251

$ws.Range("B9").Value = 27
$ws.Range("C9").Value = -1319.2560000000001
$ws.Range("D9").Value = -1338.309
$ws.Range("C9").NumberFormat = "0.00"
$ws.Range("D9").NumberFormat = "0.00"
$ws.Range("E9").Formula = "=-(C9-D9)/D9"
$ws.Range("F9").Value = 20000
$ws.Range("F9").NumberFormat = "#,##0"
$ws.Range("F9").Value = ">10,000"

# ---------------------------------------------------------------------
# New block: "Solver=SCIP" (header row 11, data rows 12-14)
# ---------------------------------------------------------------------
$ws.Range("A11").Value = "Solver=SCIP"
$ws.Range("B11").Value = "number of scenario"
$ws.Range("C11").Value = "UB"
$ws.Range("D11").Value = "LB"
$ws.Range("E11").Value = "gap"
$ws.Range("F11").Value = "Wall Time"

$ws.Range("B12").Value = 3
$ws.Range("C12").Value = -1007.7889579159799
$ws.Range("D12").Value = -2210.83461672973
$ws.Range("C12").NumberFormat = "0.00"
$ws.Range("D12").NumberFormat = "0.00"
$ws.Range("E12").Formula = "=-(C12-D12)/D12"
$ws.Range("F12").Value = ">10,000"

$ws.Range("B13").Value = 9
$ws.Range("C13").Value = -0.0000013578666711118499
$ws.Range("D13").Value = -2390.5310000300201
$ws.Range("C13").NumberFormat = "0.00"
$ws.Range("D13").NumberFormat = "0.00"
$ws.Range("E13").Formula = "=-(C13-D13)/D13"
$ws.Range("F13").Value = ">10,000"

$ws.Range("B14").Value = 27
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = -2430.51977177764
$ws.Range("C14").NumberFormat = "0.00"
$ws.Range("D14").NumberFormat = "0.00"
$ws.Range("E14").Formula = "=-(C14-D14)/D14"
$ws.Range("F14").Value = ">10,000"

# ---------------------------------------------------------------------
# The old "Number of nodes" summary table (previously rows 13-14) moves
# down to rows 21-22 to make room for the new blocks above.
# ---------------------------------------------------------------------
$ws.Range("A13:F14").ClearContents()

$ws.Range("A21").Value = "Number of nodes"
$ws.Range("B21").Value = "Stage 1 Binary"
$ws.Range("C21").Value = "Stage 2 Binary var per scenario"
$ws.Range("D21").Value = "Stage 2 continuous var per scenario"
$ws.Range("E21").Value = "Stage 2 Linear constr per scenario"
$ws.Range("F21").Value = "Nonlinear Constr per scenario"

$ws.Range("A22").Value = 9
$ws.Range("B22").Value = 9
$ws.Range("C22").Value = 35
$ws.Range("D22").Value = 112
$ws.Range("E22").Value = 116
$ws.Range("F22").Value = 22

# ---------------------------------------------------------------------
# Sheet view / selection bookkeeping: "DE" becomes the active sheet (so
# "branchAndBound" loses its tabSelected flag), with E17 selected.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("E17").Select()
